# Updates cryptocurrency price/volume data to match the latest GitHub Actions scrape.
# Numeric-looking text values are written with a leading apostrophe so Excel keeps
# them as literal text (matching the original inlineStr cell type) instead of coercing
# them to numbers and losing formatting such as trailing zeros or thousand-dot grouping.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.144.17'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.912.99'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '''324.81'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '''0.4604'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''0.3854'
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('D9').Value = '''45.67'
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').Value = '''0.07793'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('D11').Value = '''0.9688'
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').Value = '''22.27'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = '1.945.47'
$ws.Range('E13').Value = '  +3.01%  '
$ws.Range('D14').Value = '''5.752'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').Value = '''7.029'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '''0.07074'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = '''86.26'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').Value = '''1.006'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '''0.000009668'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('D20').Value = '''16.93'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').Value = '''1.003'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = '29.138.83'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = '''5.454'
$ws.Range('E23').Value = '  +2.09%  '
$ws.Range('D24').Value = '''11.05'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').Value = '2.154.27'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').Value = '''2.097'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').Value = '''157.62'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').Value = '''19.34'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').Value = '''5.705'
$ws.Range('E29').Value = '  -3.25%  '
$ws.Range('D30').Value = '''118.42'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '''1.824'
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('D32').Value = '''0.09338'
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('D33').Value = '''0.8587'
$ws.Range('E33').Value = '  -2.78%  '
$ws.Range('D34').Value = '''5.144'
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').Value = '''1.273'
$ws.Range('E35').Value = '  -3.20%  '
$ws.Range('D36').Value = '''3.077'
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('D37').Value = '''0.05758'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '''1.159'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('D39').Value = '''0.02070'
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''7.574'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.5610'
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').Value = '''0.000003070'
$ws.Range('E42').Value = '  +8.10%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '''0.1769'
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('D44').Value = '''9.290'
$ws.Range('E44').Value = '  -4.53%  '
$ws.Range('D45').Value = '''2.727'
$ws.Range('E45').Value = '  +6.84%  '
$ws.Range('D46').Value = '''0.5247'
$ws.Range('E46').Value = '  -1.51%  '
$ws.Range('D47').Value = '''11.36'
$ws.Range('E47').Value = '  -4.46%  '
$ws.Range('D48').Value = '''0.06832'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').Value = '''2.050'
$ws.Range('E49').Value = '  -6.44%  '
$ws.Range('D50').Value = '''1.796'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').Value = '''110.96'
$ws.Range('E51').Value = '  -1.58%  '

Write-Output "Applied 107 cell updates"
